# Update "想去人数" (want-to-go headcount) figures in column F across all
# four sheets of the 广州-漫展信息 workbook to the refreshed counts captured
# in this gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 2777
$ws.Range("F4").Value  = 1110
$ws.Range("F5").Value  = 20233
$ws.Range("F7").Value  = 2435
$ws.Range("F8").Value  = 766
$ws.Range("F9").Value  = 612
$ws.Range("F11").Value = 714
$ws.Range("F12").Value = 260
$ws.Range("F13").Value = 257
$ws.Range("F15").Value = 387
$ws.Range("F16").Value = 92
$ws.Range("F17").Value = 489
$ws.Range("F19").Value = 228
$ws.Range("F21").Value = 22

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value  = 38
$ws.Range("F6").Value  = 304
$ws.Range("F15").Value = 113

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6031
$ws.Range("F3").Value = 668
$ws.Range("F4").Value = 617
$ws.Range("F5").Value = 1274

# Sheet "全部类型" (All types)
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 6031
$ws.Range("F3").Value  = 668
$ws.Range("F4").Value  = 617
$ws.Range("F6").Value  = 38
$ws.Range("F7").Value  = 1274
$ws.Range("F8").Value  = 2777
$ws.Range("F9").Value  = 1110
$ws.Range("F10").Value = 20233
$ws.Range("F15").Value = 304
$ws.Range("F16").Value = 2435
$ws.Range("F17").Value = 766
$ws.Range("F19").Value = 612
$ws.Range("F21").Value = 714
$ws.Range("F22").Value = 260
$ws.Range("F23").Value = 257
$ws.Range("F28").Value = 387
$ws.Range("F29").Value = 92
$ws.Range("F32").Value = 489
$ws.Range("F36").Value = 228
$ws.Range("F37").Value = 113
$ws.Range("F38").Value = 113
$ws.Range("F43").Value = 22
